$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at row 43 (this pushes the existing weekly rows down
#    by one, new row inherits formatting from the row below per Excel's
#    default insert behaviour).
$ws.Rows.Item(43).Insert()

# Fill in the new row 43 with this week's data, copying the constant
# columns from the (now-shifted) row 44 and setting the week-specific
# values.
$ws.Cells.Item(43, 1).Value = 11
$ws.Cells.Item(43, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(43, 3).Value = "Bíobío"
$ws.Cells.Item(43, 4).Value = 45146
$ws.Cells.Item(43, 4).NumberFormat = $ws.Cells.Item(44, 4).NumberFormat
$ws.Cells.Item(43, 5).Value = 8
$ws.Cells.Item(43, 6).Value = 100114007
$ws.Cells.Item(43, 7).Value = "Jengibre"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 40
$ws.Cells.Item(43, 11).Value = 17000
$ws.Cells.Item(43, 12).Value = 18000
$ws.Cells.Item(43, 13).Value = 17500
$ws.Cells.Item(43, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(43, 15).Value = "Perú"
$ws.Cells.Item(43, 16).Value = 1346
$ws.Cells.Item(43, 17).Value = 13
$ws.Cells.Item(43, 18).Value = "Hortaliza"

# 2) Delete the (now-shifted) old row that used to be row 57 (D=44159),
#    which after the insert above sits at row 58.
$ws.Rows.Item(58).Delete()

# 3) Delete the (now-shifted) old row that used to be row 61 (D=44264),
#    which after the two edits above sits at row 61.
$ws.Rows.Item(61).Delete()
